# Filipenses integrado en archivo general
#
# Philippians (row 51) is now marked as fully incorporated ("hecho" = 1)
# instead of the provisional "*" marker, and that same "*" in-progress
# marker moves onto Amos (row 31), which had previously been "0".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Data changes -----------------------------------------------------
# Amos (book #30, row 31): hecho 0 -> "*"
$ws.Range("E31").Value = "*"

# Filipenses (book #50, row 51): hecho "*" -> 1 (done)
$ws.Range("E51").Value = 1

# --- Formatting -------------------------------------------------------
# Column E ("hecho") is right aligned.
$ws.Range("E1:E71").HorizontalAlignment = -4152

# Summary block (rows 68-71): center the count/total cells.
$ws.Range("E68").HorizontalAlignment = -4108
$ws.Range("H68").HorizontalAlignment = -4108
$ws.Range("E69").HorizontalAlignment = -4108
$ws.Range("H69").HorizontalAlignment = -4108
$ws.Range("E71").HorizontalAlignment = -4108
$ws.Range("H71").HorizontalAlignment = -4108

# Proportion row (70): bold, with the ratio cells centered too.
$ws.Range("E70:H70").Font.Bold = $true
$ws.Range("E70").HorizontalAlignment = -4108
$ws.Range("H70").HorizontalAlignment = -4108

# --- Selection / view state --------------------------------------------
$ws.Range("I72").Select()
